$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 5-6 (pushes existing rows 5+ down to 7+)
$ws.Rows("5:6").Insert()

# Populate row 6 first (Mirage) so shared-string indices land in the same
# order as the authored workbook (42/43 = Mirage, 44/45 = BroadCom),
# even though row 6 is written after row 5 is laid out.
$ws.Range("A6").Value = 'Mirage Air craft Services'
$ws.Range("B6").Value = 'E'
$ws.Range("C6").Value = '{"Result":1,"ErrTitle":"","ErrMsg":"","ErrMsgWidth":null,"ErrBtnLabel":"Ok","MultiErr":null,"AppNotice":null,"FuncName":null,"Args":[{"nameSearch":{"SearchResult":{"timestamp":"2025-09-17 14:30:18","firstName":"","lastName":"","fullName":"Mirage Air craft Services","type":"E","method":3,"tokens":{"first":[],"first_syn":[],"last":[],"last_syn":[],"full":["MIRAGE","AIR","CRAFT","SERVICES"],"full_syn":{"MIRAGE":[],"AIR":["AEREA","AEREO","AEREOS","EIR"],"CRAFT":[],"SERVICES":["SERV","SERVICIO","SERVICIOS","SERVICOS","SERVICE"]},"icij_name":["MIRAGE","AIR","CRAFT","SERVICES"],"icij_nameSyn":{"MIRAGE":[],"AIR":["AEREA","AEREO","AEREOS","EIR"],"CRAFT":[],"SERVICES":["SERV","SERVICIO","SERVICIOS","SERVICOS","SERVICE"]}},"deflated":"MIRAGEAIRCRAFTSERVICES","generated":{"pep":"2025-07-22 16:47:11","sanction":"2025-09-09 15:05:21","watch":"2025-08-07 17:33:59","mex":"2025-07-31 13:34:08","soe":"2025-06-18 10:36:28","rights":"2025-07-31 13:47:27","col":"2025-07-31 13:42:51","icij":"2016-07-15 00:00:00"},"matches":2,"hits":3,"ref":{"sanction":["361804:10131047950"],"watch":["6306178:20201576120","6444316:20201633423"]},"runtime":{"pep":0.007000000000000000145716771982051795930601656436920166015625,"sanction":0.00899999999999999931998839741709161899052560329437255859375,"watch":0.008000000000000000166533453693773481063544750213623046875,"mex":0.00899999999999999931998839741709161899052560329437255859375,"soe":0.007000000000000000145716771982051795930601656436920166015625,"rights":0.0120000000000000002498001805406602215953171253204345703125,"col":0.00899999999999999931998839741709161899052560329437255859375,"icij":0.00600000000000000012490009027033011079765856266021728515625},"error":"","countryExcl":"","listExcl":"","sqlTemplate":"SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1{filter_join}\nINNER JOIN {_tbl_} AS L2 ON L2.recid = L1.recid\nINNER JOIN {_tbl_} AS L3 ON L3.recid = L1.recid\nINNER JOIN {_tbl_} AS L4 ON L4.recid = L1.recid\nWHERE {filter_cond}L1.token = ''MIRAGE''\nAND (L2.token = ''AIR''\n  OR L2.token = ''AEREA''\n  OR L2.token = ''AEREO''\n  OR L2.token = ''AEREOS''\n  OR L2.token = ''EIR''\n)\nAND L3.token = ''CRAFT''\nAND (L4.token = ''SERVICES''\n  OR L4.token = ''SERV''\n  OR L4.token = ''SERVICIO''\n  OR L4.token = ''SERVICIOS''\n  OR L4.token = ''SERVICOS''\n  OR L4.token = ''SERVICE''\n)\nLIMIT 100","deflatedSqlTemplate":"SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1{filter_join}\nWHERE {filter_cond}L1.token = ''MIRAGEAIRCRAFTSERVICES''\nAND L1.nameSrc = ''Deflated''\n LIMIT 100","nameID":3341266,"icijName":"Mirage Air craft Services","sqlTemplate_icij":"SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1\nINNER JOIN {_tbl_} AS L2 ON L2.recid = L1.recid\nINNER JOIN {_tbl_} AS L3 ON L3.recid = L1.recid\nINNER JOIN {_tbl_} AS L4 ON L4.recid = L1.recid\nLEFT JOIN {_tbl_} AS deflateJoin ON deflateJoin.recid = L1.recid\nWHERE L1.token = ''MIRAGE''\nAND (L2.token = ''AIR''\n  OR L2.token = ''AEREA''\n  OR L2.token = ''AEREO''\n  OR L2.token = ''AEREOS''\n  OR L2.token = ''EIR''\n)\nAND L3.token = ''CRAFT''\nAND (L4.token = ''SERVICES''\n  OR L4.token = ''SERV''\n  OR L4.token = ''SERVICIO''\n  OR L4.token = ''SERVICIOS''\n  OR L4.token = ''SERVICOS''\n  OR L4.token = ''SERVICE''\n)\nAND deflateJoin.token LIKE ''%MIRAGEAIRCRAFTSERVICES%''\nLIMIT 100"},"code":"Array\n(\n    [timestamp] => 2025-09-17 14:30:18\n    [firstName] => \n    [lastName] => \n    [fullName] => Mirage Air craft Services\n    [type] => E\n    [method] => 3\n    [tokens] => Array\n        (\n            [first] => Array\n                (\n                )\n\n            [first_syn] => Array\n                (\n                )\n\n            [last] => Array\n                (\n                )\n\n            [last_syn] => Array\n                (\n                )\n\n            [full] => Array\n                (\n                    [0] => MIRAGE\n                    [1] => AIR\n                    [2] => CRAFT\n                    [3] => SERVICES\n                )\n\n            [full_syn] => Array\n                (\n                    [MIRAGE] => Array\n                        (\n                        )\n\n                    [AIR] => Array\n                        (\n                            [0] => AEREA\n                            [1] => AEREO\n                            [2] => AEREOS\n                            [3] => EIR\n                        )\n\n                    [CRAFT] => Array\n                        (\n                        )\n\n                    [SERVICES] => Array\n                        (\n                            [0] => SERV\n                            [1] => SERVICIO\n                            [2] => SERVICIOS\n                            [3] => SERVICOS\n                            [4] => SERVICE\n                        )\n\n                )\n\n            [icij_name] => Array\n                (\n                    [0] => MIRAGE\n                    [1] => AIR\n                    [2] => CRAFT\n                    [3] => SERVICES\n                )\n\n            [icij_nameSyn] => Array\n                (\n                    [MIRAGE] => Array\n                        (\n                        )\n\n                    [AIR] => Array\n                        (\n                            [0] => AEREA\n                            [1] => AEREO\n                            [2] => AEREOS\n                            [3] => EIR\n                        )\n\n                    [CRAFT] => Array\n                        (\n                        )\n\n                    [SERVICES] => Array\n                        (\n                            [0] => SERV\n                            [1] => SERVICIO\n                            [2] => SERVICIOS\n                            [3] => SERVICOS\n                            [4] => SERVICE\n                        )\n\n                )\n\n        )\n\n    [deflated] => MIRAGEAIRCRAFTSERVICES\n    [generated] => Array\n        (\n            [pep] => 2025-07-22 16:47:11\n            [sanction] => 2025-09-09 15:05:21\n            [watch] => 2025-08-07 17:33:59\n            [mex] => 2025-07-31 13:34:08\n            [soe] => 2025-06-18 10:36:28\n            [rights] => 2025-07-31 13:47:27\n            [col] => 2025-07-31 13:42:51\n            [icij] => 2016-07-15 00:00:00\n        )\n\n    [matches] => 2\n    [hits] => 3\n    [ref] => Array\n        (\n            [sanction] => Array\n                (\n                    [0] => 361804:10131047950\n                )\n\n            [watch] => Array\n                (\n                    [0] => 6306178:20201576120\n                    [1] => 6444316:20201633423\n                )\n\n        )\n\n    [runtime] => Array\n        (\n            [pep] => 0.007\n            [sanction] => 0.009\n            [watch] => 0.008\n            [mex] => 0.009\n            [soe] => 0.007\n            [rights] => 0.012\n            [col] => 0.009\n            [icij] => 0.006\n        )\n\n    [error] => \n    [countryExcl] => \n    [listExcl] => \n    [sqlTemplate] => SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1{filter_join}\nINNER JOIN {_tbl_} AS L2 ON L2.recid = L1.recid\nINNER JOIN {_tbl_} AS L3 ON L3.recid = L1.recid\nINNER JOIN {_tbl_} AS L4 ON L4.recid = L1.recid\nWHERE {filter_cond}L1.token = ''MIRAGE''\nAND (L2.token = ''AIR''\n  OR L2.token = ''AEREA''\n  OR L2.token = ''AEREO''\n  OR L2.token = ''AEREOS''\n  OR L2.token = ''EIR''\n)\nAND L3.token = ''CRAFT''\nAND (L4.token = ''SERVICES''\n  OR L4.token = ''SERV''\n  OR L4.token = ''SERVICIO''\n  OR L4.token = ''SERVICIOS''\n  OR L4.token = ''SERVICOS''\n  OR L4.token = ''SERVICE''\n)\nLIMIT 100\n    [deflatedSqlTemplate] => SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1{filter_join}\nWHERE {filter_cond}L1.token = ''MIRAGEAIRCRAFTSERVICES''\nAND L1.nameSrc = ''Deflated''\n LIMIT 100\n    [nameID] => 3341266\n    [icijName] => Mirage Air craft Services\n    [sqlTemplate_icij] => SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1\nINNER JOIN {_tbl_} AS L2 ON L2.recid = L1.recid\nINNER JOIN {_tbl_} AS L3 ON L3.recid = L1.recid\nINNER JOIN {_tbl_} AS L4 ON L4.recid = L1.recid\nLEFT JOIN {_tbl_} AS deflateJoin ON deflateJoin.recid = L1.recid\nWHERE L1.token = ''MIRAGE''\nAND (L2.token = ''AIR''\n  OR L2.token = ''AEREA''\n  OR L2.token = ''AEREO''\n  OR L2.token = ''AEREOS''\n  OR L2.token = ''EIR''\n)\nAND L3.token = ''CRAFT''\nAND (L4.token = ''SERVICES''\n  OR L4.token = ''SERV''\n  OR L4.token = ''SERVICIO''\n  OR L4.token = ''SERVICIOS''\n  OR L4.token = ''SERVICOS''\n  OR L4.token = ''SERVICE''\n)\nAND deflateJoin.token LIKE ''%MIRAGEAIRCRAFTSERVICES%''\nLIMIT 100\n)\n","Preview":{"pep":[],"sanction":[{"recid":361804,"ID":"10131047950","First_Name":"","Last_Name":"MIRAGE AIR CRAFT SERVICES SOLE PROPRIETORSHIP LLC","Full_Name":"","Other_Names":"","AltScript":"","RecType":"E"}],"watch":[{"recid":6306178,"ID":"20201576120","First_Name":"","Last_Name":"","Full_Name":"MIRAGE AIR CRAFT SERVICES SOLE PROPREITORSHIP LLC","Other_Names":"","AltScript":"","RecType":"E"},{"recid":6444316,"ID":"20201633423","First_Name":"","Last_Name":"","Full_Name":"MIRAGE AIR CRAFT SERVICES SOLE PROPRIETORSHIP LLC","Other_Names":"","AltScript":"","RecType":"E"}],"icij":[],"mex":[],"soe":[],"rights":[],"col":[],"media":[],"ofac":[]}}}],"Redirect":null}'

# Row 5: BroadCom Wireless Communications Corporation (bold name, new style)
$ws.Range("A5").ClearFormats()
$ws.Range("A5").Value = 'BroadCom Wireless Communications Corporation'
$ws.Range("A5").Font.Bold = $true
$ws.Range("B5").Value = 'E'
$ws.Range("C5").Value = '{"Result":1,"ErrTitle":"","ErrMsg":"","ErrMsgWidth":null,"ErrBtnLabel":"Ok","MultiErr":null,"AppNotice":null,"FuncName":null,"Args":[{"nameSearch":{"SearchResult":{"timestamp":"2025-09-17 14:35:05","firstName":"","lastName":"","fullName":"BroadCom Wireless Communications Corporation","type":"E","method":3,"tokens":{"first":[],"first_syn":[],"last":[],"last_syn":[],"full":["BROADCOM","WIRELESS","COMMUNICATIONS"],"full_syn":{"BROADCOM":[],"WIRELESS":[],"COMMUNICATIONS":[]},"icij_name":["BROADCOM","WIRELESS","COMMUNICATIONS"],"icij_nameSyn":{"BROADCOM":[],"WIRELESS":[],"COMMUNICATIONS":[]}},"deflated":"BROADCOMWIRELESSCOMMUNICATIONSCORPORATION","generated":{"pep":"2025-07-22 16:47:11","sanction":"2025-09-09 15:05:21","watch":"2025-08-07 17:33:59","mex":"2025-07-31 13:34:08","soe":"2025-06-18 10:36:28","rights":"2025-07-31 13:47:27","col":"2025-07-31 13:42:51","icij":"2016-07-15 00:00:00"},"matches":1,"hits":1,"ref":{"watch":["3319709:20200616696"]},"runtime":{"pep":0.01400000000000000029143354396410359186120331287384033203125,"sanction":0.0290000000000000014710455076283324160613119602203369140625,"watch":0.0259999999999999988065102485279567190445959568023681640625,"mex":0.0189999999999999995281552145343084703199565410614013671875,"soe":0.0259999999999999988065102485279567190445959568023681640625,"rights":0.00899999999999999931998839741709161899052560329437255859375,"col":0.007000000000000000145716771982051795930601656436920166015625,"icij":0.01000000000000000020816681711721685132943093776702880859375},"error":"","countryExcl":"","listExcl":"","sqlTemplate":"SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1{filter_join}\nINNER JOIN {_tbl_} AS L2 ON L2.recid = L1.recid\nINNER JOIN {_tbl_} AS L3 ON L3.recid = L1.recid\nWHERE {filter_cond}L1.token = ''BROADCOM''\nAND L2.token = ''WIRELESS''\nAND L3.token = ''COMMUNICATIONS''\nLIMIT 100","deflatedSqlTemplate":"SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1{filter_join}\nWHERE {filter_cond}L1.token = ''BROADCOMWIRELESSCOMMUNICATIONSCORPORATION''\nAND L1.nameSrc = ''Deflated''\n LIMIT 100","nameID":3556807,"icijName":"BroadCom Wireless Communications Corporation","sqlTemplate_icij":"SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1\nINNER JOIN {_tbl_} AS L2 ON L2.recid = L1.recid\nINNER JOIN {_tbl_} AS L3 ON L3.recid = L1.recid\nLEFT JOIN {_tbl_} AS deflateJoin ON deflateJoin.recid = L1.recid\nWHERE L1.token = ''BROADCOM''\nAND L2.token = ''WIRELESS''\nAND L3.token = ''COMMUNICATIONS''\nAND deflateJoin.token LIKE ''%BROADCOMWIRELESSCOMMUNICATIONSCORPORATION%''\nLIMIT 100"},"code":"Array\n(\n    [timestamp] => 2025-09-17 14:35:05\n    [firstName] => \n    [lastName] => \n    [fullName] => BroadCom Wireless Communications Corporation\n    [type] => E\n    [method] => 3\n    [tokens] => Array\n        (\n            [first] => Array\n                (\n                )\n\n            [first_syn] => Array\n                (\n                )\n\n            [last] => Array\n                (\n                )\n\n            [last_syn] => Array\n                (\n                )\n\n            [full] => Array\n                (\n                    [0] => BROADCOM\n                    [1] => WIRELESS\n                    [2] => COMMUNICATIONS\n                )\n\n            [full_syn] => Array\n                (\n                    [BROADCOM] => Array\n                        (\n                        )\n\n                    [WIRELESS] => Array\n                        (\n                        )\n\n                    [COMMUNICATIONS] => Array\n                        (\n                        )\n\n                )\n\n            [icij_name] => Array\n                (\n                    [0] => BROADCOM\n                    [1] => WIRELESS\n                    [2] => COMMUNICATIONS\n                )\n\n            [icij_nameSyn] => Array\n                (\n                    [BROADCOM] => Array\n                        (\n                        )\n\n                    [WIRELESS] => Array\n                        (\n                        )\n\n                    [COMMUNICATIONS] => Array\n                        (\n                        )\n\n                )\n\n        )\n\n    [deflated] => BROADCOMWIRELESSCOMMUNICATIONSCORPORATION\n    [generated] => Array\n        (\n            [pep] => 2025-07-22 16:47:11\n            [sanction] => 2025-09-09 15:05:21\n            [watch] => 2025-08-07 17:33:59\n            [mex] => 2025-07-31 13:34:08\n            [soe] => 2025-06-18 10:36:28\n            [rights] => 2025-07-31 13:47:27\n            [col] => 2025-07-31 13:42:51\n            [icij] => 2016-07-15 00:00:00\n        )\n\n    [matches] => 1\n    [hits] => 1\n    [ref] => Array\n        (\n            [watch] => Array\n                (\n                    [0] => 3319709:20200616696\n                )\n\n        )\n\n    [runtime] => Array\n        (\n            [pep] => 0.014\n            [sanction] => 0.029\n            [watch] => 0.026\n            [mex] => 0.019\n            [soe] => 0.026\n            [rights] => 0.009\n            [col] => 0.007\n            [icij] => 0.01\n        )\n\n    [error] => \n    [countryExcl] => \n    [listExcl] => \n    [sqlTemplate] => SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1{filter_join}\nINNER JOIN {_tbl_} AS L2 ON L2.recid = L1.recid\nINNER JOIN {_tbl_} AS L3 ON L3.recid = L1.recid\nWHERE {filter_cond}L1.token = ''BROADCOM''\nAND L2.token = ''WIRELESS''\nAND L3.token = ''COMMUNICATIONS''\nLIMIT 100\n    [deflatedSqlTemplate] => SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1{filter_join}\nWHERE {filter_cond}L1.token = ''BROADCOMWIRELESSCOMMUNICATIONSCORPORATION''\nAND L1.nameSrc = ''Deflated''\n LIMIT 100\n    [nameID] => 3556807\n    [icijName] => BroadCom Wireless Communications Corporation\n    [sqlTemplate_icij] => SELECT DISTINCT L1.recid\nFROM {_tbl_} AS L1\nINNER JOIN {_tbl_} AS L2 ON L2.recid = L1.recid\nINNER JOIN {_tbl_} AS L3 ON L3.recid = L1.recid\nLEFT JOIN {_tbl_} AS deflateJoin ON deflateJoin.recid = L1.recid\nWHERE L1.token = ''BROADCOM''\nAND L2.token = ''WIRELESS''\nAND L3.token = ''COMMUNICATIONS''\nAND deflateJoin.token LIKE ''%BROADCOMWIRELESSCOMMUNICATIONSCORPORATION%''\nLIMIT 100\n)\n","Preview":{"pep":[],"sanction":[],"watch":[{"recid":3319709,"ID":"20200616696","First_Name":"","Last_Name":"","Full_Name":"BroadCom Wireless Communications Corporation","Other_Names":"","AltScript":"","RecType":"E"}],"icij":[],"mex":[],"soe":[],"rights":[],"col":[],"media":[],"ofac":[]}}}],"Redirect":null}'

# Widen column C to fit the new content
$ws.Columns("C").ColumnWidth = 93.5 - 5/6

# Update the active selection to match the final cursor position
$ws.Range("C22").Select()

